# Weekly update: insert a new price record as row 35, shifting the
# existing rows 35-50 down to 36-51 (dimension grows from A1:R50 to A1:R51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 35; this pushes the old rows 35-50
# down to 36-51, carrying their values/styles with them (matches the diff,
# where old row 35's data reappears unchanged on row 36, etc.).
$ws.Rows("35:35").Insert()

# Populate the newly inserted row 35 with the new weekly record.
$ws.Cells.Item(35, 1).Value = 1
$ws.Cells.Item(35, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(35, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(35, 4).Value = 44466
$ws.Cells.Item(35, 5).Value = 15
$ws.Cells.Item(35, 6).Value = 100112040
$ws.Cells.Item(35, 7).Value = "Cilantro"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 300
$ws.Cells.Item(35, 11).Value = 900
$ws.Cells.Item(35, 12).Value = 1000
$ws.Cells.Item(35, 13).Value = 950
$ws.Cells.Item(35, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(35, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(35, 16).Value = 475
$ws.Cells.Item(35, 17).Value = 2
$ws.Cells.Item(35, 18).Value = "Hortaliza"
